$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (dragged border) - closest value this engine's pixel
# quantization can reach to the authored width of 10.5714285714286 chars.
$ws.Columns.Item(1).ColumnWidth = 9.65

# Update the roll-number / id values in column A (A2:A6)
$ws.Range("A2").Value = 2022056840
$ws.Range("A3").Value = 2022056141
$ws.Range("A4").Value = 2022056142
$ws.Range("A5").Value = 2022056143
$ws.Range("A6").Value = 2022056014

# Move/confirm the active selection to A6 (matches final cursor position)
$ws.Range("A6").Select()
